# Insert a new column before D (shifts D..K -> E..L), carrying over
# number formats from the (new) column E, then populate column D with the
# new reporting-period figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a full column at D; everything from D..K shifts right to E..L.
$ws.Range("D1").EntireColumn.Insert()

# Copy the number/date formatting from column E (the old column D, now
# shifted one to the right) into the freshly inserted column D so the new
# cells pick up the same styles (date format row 7/38/80, number format
# elsewhere) instead of the generic default style.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newValues = @{
    7 = 43465
    8 = 152000
    9 = 3700
    10 = 148300
    13 = 0
    14 = 0
    15 = 400
    17 = 115700
    18 = 36300
    20 = 0
    21 = 47600
    22 = 0
    23 = 36300
    24 = -67900
    25 = 0
    26 = 104300
    27 = 104300
    28 = 0
    29 = -74000
    30 = 0
    31 = 0
    32 = 0
    33 = 30300
    34 = 0
    35 = 30300
    38 = 43465
    41 = 194300
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 1189500
    48 = 22000
    49 = 14100
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 1695200
    57 = 694400
    58 = 408700
    59 = 313300
    60 = 0
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 1465300
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 75200
    73 = 0
    74 = 0
    75 = 0
    76 = 229800
    77 = 0
    80 = 43465
    81 = 30300
    83 = 11200
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 48300
    91 = -4400
    92 = 0
    93 = 0
    94 = -46800
    96 = -13200
    97 = 0
    98 = 0
    99 = 0
    100 = 1300
    101 = 0
    102 = 2800
    12 = "NA"
}


foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}
